$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G8").Value = 1.45
$ws.Range("L8").Value = 6.5
$ws.Range("AA8").Value = 13
$ws.Range("AD8").Value = 8
$ws.Range("AE8").Value = 19
$ws.Range("AX8").Value = 34
$ws.Range("BB8").Value = 301
